# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2-63) holds a date serial meant to represent the
# "as-of" snapshot date for each forecast row. Every stored value was
# off by a month (it landed on the 1st of the reporting month instead
# of the 15th of the following month); this re-points each date to the
# 15th of the month after the one currently stored.
#
# We lean on Excel's own DATE/YEAR/MONTH functions via
# Application.Evaluate (not a helper cell + Range.Formula) so no stray
# cell/style is left behind in the workbook - Evaluate just computes a
# value without touching any cell's content or number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2

    $new = $excel.Evaluate("DATE(YEAR($old),MONTH($old)+1,15)")

    $cell.Value2 = $new
}
